$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2 = @{ B = 3.286832544864788;  C = 1.655778082260271;  D = 0.1494219747398047; E = 0.4942365360607697; F = 0; G = 5.586269137925634 }
    3 = @{ B = 0.1190320826869504; C = 0.306821227259698;   D = 0.1494219747398047; E = 0.4942365360607697; F = 0; G = 1.069511820747223 }
    4 = @{ B = 3.286832544864788;  C = 1.655778082260271;  D = 0.7527432677738641; E = 0.4942365360607697; F = 1; G = 6.189590430959694 }
    5 = @{ B = 0.003208871385164791; C = 0.306821227259698; D = 0.7527432677738641; E = 10.19245300693656;  F = 0; G = 11.25522637335528 }
    6 = @{ B = 0.01293466051926884;  C = 0.306821227259698; D = 0.1494219747398047; E = 0.4942365360607697; F = 1; G = 0.9634143985795411 }
    7 = @{ B = 3.286832544864788;  C = 1.655778082260271;  D = 0.7527432677738641; E = 0.4942365360607697; F = 0; G = 6.189590430959694 }
    8 = @{ B = 0.6606524410359556; C = 10.34677158129881;  D = 22.3905356188092;   E = 10.19245300693656;  F = 1; G = 43.59041264808052 }
    9 = @{ B = 3.286832544864788;  C = 1.655778082260271;  D = 0.1494219747398047; E = 0.4942365360607697; F = 1; G = 5.586269137925634 }
}

foreach ($row in $values.Keys) {
    $rowVals = $values[$row]
    $ws.Range("B$row").Value = $rowVals.B
    $ws.Range("C$row").Value = $rowVals.C
    $ws.Range("D$row").Value = $rowVals.D
    $ws.Range("E$row").Value = $rowVals.E
    $ws.Range("F$row").Value = $rowVals.F
    $ws.Range("G$row").Value = $rowVals.G
}
